$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.646.26'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '1.563.19'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Value = '''210.60'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').Value = '''0.489'
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('D8').Value = '''25.15'
$ws.Range('E8').Value = '  +5.53%  '
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '1.785.89'
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('D13').Value = '1.574.77'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').Value = '28.629.83'
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').Value = '''0.515'
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').Value = '''3.64'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').Value = '''61.27'
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').Value = '''229.52'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('E20').Value = '  -1.45%  '
$ws.Range('D21').Value = '''1.00'
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('E22').Value = '  -1.16%  '
$ws.Range('D23').Value = '''9.00'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').Value = '''150.90'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = '''14.78'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('D29').Value = '''6.23'
$ws.Range('E29').Value = '  -2.20%  '
$ws.Range('E30').Value = '  -4.19%  '
$ws.Range('E31').Value = '  -2.48%  '
$ws.Range('E32').Value = '  -0.53%  '
$ws.Range('D33').Value = '1.387.34'
$ws.Range('E33').Value = '  +0.34%  '
$ws.Range('D34').Value = '''2.98'
$ws.Range('E34').Value = '  -4.28%  '
$ws.Range('D35').Value = '''1.02'
$ws.Range('E35').Value = '  -4.39%  '
$ws.Range('E36').Value = '  -1.47%  '
$ws.Range('D37').Value = '''2.68'
$ws.Range('E37').Value = '  +0.80%  '
$ws.Range('E38').Value = '  -2.35%  '
$ws.Range('E39').Value = '  -1.32%  '
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('D41').Value = '''0.518'
$ws.Range('E41').Value = '  -0.96%  '
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').Value = '''0.772'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('E44').Value = '  -2.27%  '
$ws.Range('D45').Value = '''64.02'
$ws.Range('E45').Value = '  +2.47%  '
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('D47').Value = '1.698.88'
$ws.Range('E47').Value = '  -0.69%  '
$ws.Range('D48').Value = '''0.870'
$ws.Range('E48').Value = '  -5.43%  '
$ws.Range('D49').Value = '''85.11'
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('D50').Value = '''43.26'
$ws.Range('E50').Value = '  +6.57%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.0512'
$ws.Range('E51').Value = '  -0.42%  '
